# Fixed errors with MSM implementation.
# Target sheet: "methodNumberOfLines" (11th sheet) holds one row per method
# with its line count in column C. This edit:
#   1. Removes the seven constructor rows entirely.
#   2. Corrects the line count from 1 to 0 for compiler-generated /
#      single-expression-lambda methods that were being mis-counted.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("methodNumberOfLines")

# --- Step 1: materialise a clean shared-string "0" cell we can replicate from.
# A plain `.Value = "0"` is auto-coerced to a numeric cell, so force Text
# format first, write the text, then clear the (now superfluous) formatting
# back off the cell -- it stays a text cell (shared string) afterwards.
$helper = $ws.Range("Z1")
$helper.NumberFormat = "@"
$helper.Value = "0"
$helper.ClearFormats()

# --- Step 2: for every row whose "Number of Lines" (column C) is the text
# "1" but is NOT one of the no-arg-constructor rows we are about to delete
# outright, replace it with the text "0" (copied from the helper cell so the
# shared-string/text typing - not a numeric type - is preserved).
$rowsToZero = @(2, 3, 6, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 29, 30, 36, 37, 38)
foreach ($r in $rowsToZero) {
    $helper.Copy($ws.Cells.Item($r, 3))
}

# --- Step 3: remove the helper cell's content so it doesn't linger in the sheet.
$helper.ClearContents()

# --- Step 4: delete the rows that correspond to no-arg/explicit constructors
# (KeycloakRole, GatewayConfiguration, CrossOriginRequestSharingFilter,
# GatewayApplication, GatewayController, KeycloakRealmRoleConverter,
# SecurityConfiguration). Delete bottom-to-top so earlier row numbers stay valid.
$rowsToDelete = @(35, 27, 25, 23, 21, 8, 4)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}
